$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2176.4666
$ws.Range("I39").Value = 2617.3333
$ws.Range("K39").Value = 7851.999899999999
$ws.Range("M39").Value = -7555.999899999999

$ws.Range("H52").Value = 5000
$ws.Range("J52").Value = 5000
$ws.Range("L52").Value = 15000
$ws.Range("N52").Value = -15320

$ws.Range("H101").Value = 546.8
$ws.Range("I101").Value = 546.8
$ws.Range("K101").Value = 1640.4
$ws.Range("M101").Value = -18.39999999999986

$ws.Range("H137").Value = 1570.7333
$ws.Range("I137").Value = 1456.2
$ws.Range("J137").Value = 1799.8
$ws.Range("K137").Value = 4368.6
$ws.Range("L137").Value = 5399.4
$ws.Range("M137").Value = -1818.6
$ws.Range("N137").Value = -10499.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1770.2632
$ws.Range("I2").Value = 653.1429000000001
$ws.Range("J2").Value = 4898.2
$ws.Range("K2").Value = 653.1429000000001
$ws.Range("L2").Value = 4898.2
$ws.Range("M2").Value = -540.1429000000001
$ws.Range("N2").Value = -5124.2

$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H74").Value = 1099.75
$ws.Range("I74").Value = 1099.75
$ws.Range("K74").Value = 1099.75
$ws.Range("M74").Value = -225.75

$ws.Range("H77").Value = 1099.75
$ws.Range("I77").Value = 1099.75
$ws.Range("K77").Value = 5498.75
$ws.Range("M77").Value = -1130.75

$ws.Range("H95").Value = 34444
$ws.Range("J95").Value = 34444
$ws.Range("L95").Value = 34444
$ws.Range("N95").Value = -39936

$ws.Range("H97").Value = 1111.1111
$ws.Range("I97").Value = 1100
$ws.Range("K97").Value = 1100
$ws.Range("M97").Value = -604

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H116").Value = 1770.2632
$ws.Range("I116").Value = 653.1429000000001
$ws.Range("J116").Value = 4898.2
$ws.Range("K116").Value = 653.1429000000001
$ws.Range("L116").Value = 4898.2
$ws.Range("M116").Value = 1640.8571
$ws.Range("N116").Value = -9486.200000000001

$ws.Range("H132").Value = 1440
$ws.Range("I132").Value = 1467.2
$ws.Range("K132").Value = 4401.6
$ws.Range("M132").Value = -1871.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1770.2632
$ws.Range("I3").Value = 653.1429000000001
$ws.Range("J3").Value = 4898.2
$ws.Range("K3").Value = 653.1429000000001
$ws.Range("L3").Value = 4898.2
$ws.Range("M3").Value = -539.1429000000001
$ws.Range("N3").Value = -5126.2

$ws.Range("H105").Value = 2293.4285
$ws.Range("I105").Value = 2119.5625
$ws.Range("J105").Value = 2849.8
$ws.Range("K105").Value = 2119.5625
$ws.Range("L105").Value = 2849.8
$ws.Range("M105").Value = -372.5625
$ws.Range("N105").Value = -6343.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5777.6665
$ws.Range("J4").Value = 7199.8
$ws.Range("L4").Value = 7199.8
$ws.Range("N4").Value = -7423.8

$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2609

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2818

$ws.Range("H88").Value = 37512.09
$ws.Range("J88").Value = 37512.09
$ws.Range("L88").Value = 37512.09
$ws.Range("N88").Value = -38324.09

$ws.Range("H91").Value = 37512.09
$ws.Range("J91").Value = 37512.09
$ws.Range("L91").Value = 37512.09
$ws.Range("N91").Value = -40320.09

$ws.Range("H92").Value = 48000
$ws.Range("J92").Value = 48000
$ws.Range("L92").Value = 48000
$ws.Range("N92").Value = -52992

$ws.Range("H122").Value = 3104.8823
$ws.Range("I122").Value = 2673.9375
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 8021.8125
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -5571.8125
$ws.Range("N122").Value = -34900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2138.2
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 2650.25
$ws.Range("K17").Value = 270
$ws.Range("L17").Value = 7950.75
$ws.Range("M17").Value = -101
$ws.Range("N17").Value = -8288.75

$ws.Range("H69").Value = 3033.1667
$ws.Range("J69").Value = 3579.8
$ws.Range("L69").Value = 10739.4
$ws.Range("N69").Value = -12361.4

$ws.Range("H72").Value = 3033.1667
$ws.Range("J72").Value = 3579.8
$ws.Range("L72").Value = 32218.2
$ws.Range("N72").Value = -40330.2

$ws.Range("H129").Value = 1272.75
$ws.Range("I129").Value = 1030.3334
$ws.Range("K129").Value = 3091.0002
$ws.Range("M129").Value = 1908.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 12863.4
$ws.Range("J92").Value = 9829.25
$ws.Range("L92").Value = 9829.25
$ws.Range("N92").Value = -13573.25

$ws.Range("H122").Value = 36297.93
$ws.Range("I122").Value = 1278
$ws.Range("K122").Value = 3834
$ws.Range("M122").Value = -1384

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1414.8889
$ws.Range("I16").Value = 1333.3636
$ws.Range("J16").Value = 1543
$ws.Range("K16").Value = 1333.3636
$ws.Range("L16").Value = 1543
$ws.Range("M16").Value = -1163.3636
$ws.Range("N16").Value = -1883

$ws.Range("H46").Value = 999
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811

$ws.Range("H68").Value = 3401
$ws.Range("I68").Value = 2600
$ws.Range("K68").Value = 2600
$ws.Range("M68").Value = -1851

$ws.Range("H71").Value = 3401
$ws.Range("I71").Value = 2600
$ws.Range("K71").Value = 13000
$ws.Range("M71").Value = -9256

$ws.Range("H76").Value = 14557.4
$ws.Range("J76").Value = 14557.4
$ws.Range("L76").Value = 14557.4
$ws.Range("N76").Value = -15233.4

$ws.Range("H79").Value = 14557.4
$ws.Range("J79").Value = 14557.4
$ws.Range("L79").Value = 14557.4
$ws.Range("N79").Value = -16897.4

$ws.Range("H93").Value = 1385.35
$ws.Range("I93").Value = 1260.8667
$ws.Range("J93").Value = 1758.8
$ws.Range("K93").Value = 1260.8667
$ws.Range("L93").Value = 1758.8
$ws.Range("M93").Value = -12.86670000000004
$ws.Range("N93").Value = -4254.8

$ws.Range("H105").Value = 63999
$ws.Range("J105").Value = 63999
$ws.Range("L105").Value = 63999
$ws.Range("N105").Value = -70987

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H107").Value = 1700
$ws.Range("J107").Value = 1549.5
$ws.Range("L107").Value = 4648.5
$ws.Range("N107").Value = -8488.5

$ws.Range("H113").Value = 967.875
$ws.Range("I113").Value = 973.75
$ws.Range("J113").Value = 962
$ws.Range("K113").Value = 2921.25
$ws.Range("L113").Value = 2886
$ws.Range("M113").Value = -751.25
$ws.Range("N113").Value = -7226

$ws.Range("H126").Value = 3262.125
$ws.Range("J126").Value = 4300
$ws.Range("L126").Value = 12900
$ws.Range("N126").Value = -17840
